# Applies the cryptocurrency price/volume update described by the commit
# "Updated cryptos list ... with GitHub Actions".
#
# All data cells in columns B:E are stored as plain text in the workbook
# (prices/percentages are formatted strings, not numbers). Assigning a
# number-looking string straight to Range.Value lets Excel auto-convert it
# to a real number (e.g. "237.30" -> 237.3), which would corrupt the data.
# Prefixing the text with a literal apostrophe forces Excel to keep it as
# text (same effect as typing `'value` into a cell); resetting the style
# back to "Normal" afterwards clears the "number stored as text" marker so
# the cell formatting matches the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $text) {
    $ws.Range($address).Value = "'" + $text
    $ws.Range($address).Style = "Normal"
}

Set-TextValue "D2" "93.433.56"
Set-TextValue "E2" "  +2.21%  "

Set-TextValue "D3" "3.126.19"
Set-TextValue "E3" "  +0.73%  "

Set-TextValue "E4" "  -0.05%  "

Set-TextValue "D5" "237.30"
Set-TextValue "E5" "  -2.45%  "

Set-TextValue "D6" "612.88"
Set-TextValue "E6" "  -0.38%  "

Set-TextValue "D7" "1.11"
Set-TextValue "E7" "  +1.45%  "

Set-TextValue "D8" "0.393"
Set-TextValue "E8" "  +2.76%  "

Set-TextValue "D9" "1.00"
Set-TextValue "E9" "  -0.10%  "

Set-TextValue "D10" "0.838"
Set-TextValue "E10" "  +13.08%  "

Set-TextValue "D11" "3.127.52"
Set-TextValue "E11" "  +0.77%  "

Set-TextValue "E12" "  -2.86%  "

Set-TextValue "D13" "0.0000246"
Set-TextValue "E13" "  -1.06%  "

Set-TextValue "D14" "35.10"
Set-TextValue "E14" "  +1.77%  "

Set-TextValue "D15" "93.182.56"
Set-TextValue "E15" "  +1.86%  "

Set-TextValue "D16" "5.45"
Set-TextValue "E16" "  -2.74%  "

Set-TextValue "D17" "3.708.43"
Set-TextValue "E17" "  +0.68%  "

Set-TextValue "D18" "3.112.42"
Set-TextValue "E18" "  -0.71%  "

Set-TextValue "D19" "3.76"
Set-TextValue "E19" "  +0.95%  "

Set-TextValue "D20" "14.73"
Set-TextValue "E20" "  -0.02%  "

Set-TextValue "D21" "5.98"
Set-TextValue "E21" "  +3.85%  "

Set-TextValue "E22" "  +1.08%  "

Set-TextValue "D23" "443.23"
Set-TextValue "E23" "  -0.65%  "

Set-TextValue "D24" "9.18"
Set-TextValue "E24" "  -0.80%  "

Set-TextValue "D25" "5.70"
Set-TextValue "E25" "  +1.44%  "

Set-TextValue "D26" "12.78"
Set-TextValue "E26" "  +9.96%  "

Set-TextValue "D27" "86.20"
Set-TextValue "E27" "  -3.00%  "

Set-TextValue "D28" "3.295.06"
Set-TextValue "E28" "  +0.13%  "

Set-TextValue "D29" "0.999"
Set-TextValue "E29" "  -0.04%  "

Set-TextValue "D30" "0.181"
Set-TextValue "E30" "  +8.53%  "

Set-TextValue "D31" "0.237"
Set-TextValue "E31" "  +4.39%  "

Set-TextValue "E32" "  -12.61%  "

Set-TextValue "E33" "  +4.35%  "

Set-TextValue "D34" "9.23"
Set-TextValue "E34" "  -0.07%  "

Set-TextValue "D35" "8.13"
Set-TextValue "E35" "  +6.20%  "

Set-TextValue "E36" "  -8.42%  "

Set-TextValue "D37" "25.98"
Set-TextValue "E37" "  -0.55%  "

Set-TextValue "D38" "3.97"
Set-TextValue "E38" "  +0.73%  "

Set-TextValue "E39" "  -1.43%  "

Set-TextValue "E40" "  +0.54%  "

Set-TextValue "D41" "24.03"
Set-TextValue "E41" "  +8.20%  "

Set-TextValue "D42" "475.20"
Set-TextValue "E42" "  -2.31%  "

Set-TextValue "D43" "0.439"
Set-TextValue "E43" "  +1.63%  "

Set-TextValue "D44" "3.33"
Set-TextValue "E44" "  -2.26%  "

Set-TextValue "E45" "  +0.00%  "

Set-TextValue "D46" "159.28"
Set-TextValue "E46" "  -0.17%  "

Set-TextValue "D47" "0.692"
Set-TextValue "E47" "  -0.23%  "

Set-TextValue "E48" "  -2.09%  "

Set-TextValue "D49" "1.34"
Set-TextValue "E49" "  +0.47%  "

Set-TextValue "B50" "Filecoin"
Set-TextValue "C50" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D50" "4.44"
Set-TextValue "E50" "  +1.66%  "

Set-TextValue "B51" "OKB"
Set-TextValue "C51" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D51" "44.04"
Set-TextValue "E51" "  -0.10%  "
